$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
# Apply as text to preserve formatting like "1.000" / "26.141.26" without numeric coercion.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.141.26"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.748.52"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.13"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5292"
$ws.Range("E7").Value = "  +1.73%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2807"
$ws.Range("E8").Value = "  -1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06178"
$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.748.29"
$ws.Range("E10").Value = "  -0.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07180"
$ws.Range("E11").Value = "  +2.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.46"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6460"
$ws.Range("E13").Value = "  +0.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.632"
$ws.Range("E14").Value = "  +2.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.41"
$ws.Range("E15").Value = "  +1.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.032.15"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.78"
$ws.Range("E19").Value = "  +2.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006765"
$ws.Range("E20").Value = "  +2.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.970.82"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.336"
$ws.Range("E22").Value = "  +4.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.737"
$ws.Range("E23").Value = "  +0.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.233"
$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.59"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.523"
$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.27"
$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.805"
$ws.Range("E28").Value = "  -2.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.03"
$ws.Range("E29").Value = "  +1.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08302"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.772"
$ws.Range("E31").Value = "  +3.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.633"
$ws.Range("E32").Value = "  +5.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04631"
$ws.Range("E33").Value = "  +4.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.643"
$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.006"
$ws.Range("E35").Value = "  +1.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6313"
$ws.Range("E36").Value = "  +3.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("E37").Value = "  +0.66%  "

$ws.Range("E38").Value = "  +2.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.979"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9991"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.88"
$ws.Range("E41").Value = "  +1.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3924"
$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7526"
$ws.Range("E43").Value = "  +2.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.079"
$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1152"
$ws.Range("E45").Value = "  +3.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.345"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05342"
$ws.Range("E47").Value = "  -2.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.62"
$ws.Range("E48").Value = "  +3.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.05"
$ws.Range("E49").Value = "  +3.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3483"
$ws.Range("E50").Value = "  +1.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.606"
$ws.Range("E51").Value = "  +0.30%  "
